$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 7 (pushes existing rows 7.. down by one)
$ws.Rows("7:7").Insert()

# Populate the new row 7 with the cloudAppUrl config entry
$ws.Range("A7").Value = "*"
$ws.Range("B7").Value = "general"
$ws.Range("C7").Value = "cloudAppUrl"
$ws.Range("D7").Formula = '="@@."&A7&"."&B7&"."&C7&"@@"'
$ws.Range("E7").Value = "Test"
$ws.Range("F7").Value = "Test"
$ws.Range("G7").Value = "Test"
$ws.Range("H7").Value = "tes"
$ws.Range("I7").Value = ""
